$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2928123333333333
$ws.Range("H2").Value = 0.8784369999999999
$ws.Range("I2").Value = 0.06406943071632207
$ws.Range("J2").Value = 0.06406943071632207
$ws.Range("M2").Value = 14.440165
$ws.Range("N2").Value = 43.320495
$ws.Range("O2").Value = 0.1441015470002482
$ws.Range("P2").Value = 0.1441015470002482
$ws.Range("Q2").Value = 4.228258407368333
$ws.Range("R2").Value = 38.05432566631499
$ws.Range("S2").Value = 0.009232504081647233
$ws.Range("T2").Value = 0.009232504081647231

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2928123333333333
$ws.Range("H3").Value = 0.8784369999999999
$ws.Range("I3").Value = 0.06406943071632207
$ws.Range("J3").Value = 0.06406943071632207
$ws.Range("O3").Value = 0.3846359116098663
$ws.Range("P3").Value = 0.3846359116098662
$ws.Range("Q3").Value = 11.28606917063422
$ws.Range("R3").Value = 101.574622535708
$ws.Range("S3").Value = 0.02464340388989771
$ws.Range("T3").Value = 0.0246434038898977

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2928123333333333
$ws.Range("H4").Value = 0.8784369999999999
$ws.Range("I4").Value = 0.06406943071632207
$ws.Range("J4").Value = 0.06406943071632207
$ws.Range("M4").Value = 21.954262
$ws.Range("N4").Value = 65.862786
$ws.Range("O4").Value = 0.2190863551385157
$ws.Range("P4").Value = 0.2190863551385156
$ws.Range("Q4").Value = 6.428478682831332
$ws.Range("R4").Value = 57.856308145482
$ws.Range("S4").Value = 0.01403673805143866
$ws.Range("T4").Value = 0.01403673805143866

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2928123333333333
$ws.Range("H5").Value = 0.8784369999999999
$ws.Range("I5").Value = 0.06406943071632207
$ws.Range("J5").Value = 0.06406943071632207
$ws.Range("M5").Value = 25.27013633333333
$ws.Range("N5").Value = 75.81040899999999
$ws.Range("O5").Value = 0.2521761862513699
$ws.Range("P5").Value = 0.2521761862513699
$ws.Range("Q5").Value = 7.399407583414775
$ws.Range("R5").Value = 66.59466825073299
$ws.Range("S5").Value = 0.01615678469333848
$ws.Range("T5").Value = 0.01615678469333848

# Row 6
$ws.Range("I6").Value = 0.8630927339690215
$ws.Range("J6").Value = 0.8630927339690215
$ws.Range("M6").Value = 14.440165
$ws.Range("N6").Value = 43.320495
$ws.Range("O6").Value = 0.1441015470002482
$ws.Range("P6").Value = 0.1441015470002482
$ws.Range("Q6").Value = 56.95975550182833
$ws.Range("R6").Value = 512.6377995164549
$ws.Range("S6").Value = 0.1243729981696097
$ws.Range("T6").Value = 0.1243729981696097

# Row 7
$ws.Range("I7").Value = 0.8630927339690215
$ws.Range("J7").Value = 0.8630927339690215
$ws.Range("O7").Value = 0.3846359116098663
$ws.Range("P7").Value = 0.3846359116098662
$ws.Range("S7").Value = 0.3319764605340264
$ws.Range("T7").Value = 0.3319764605340264

# Row 8
$ws.Range("I8").Value = 0.8630927339690215
$ws.Range("J8").Value = 0.8630927339690215
$ws.Range("M8").Value = 21.954262
$ws.Range("N8").Value = 65.862786
$ws.Range("O8").Value = 0.2190863551385157
$ws.Range("P8").Value = 0.2190863551385156
$ws.Range("Q8").Value = 86.59938413051933
$ws.Range("R8").Value = 779.3944571746739
$ws.Range("S8").Value = 0.1890918412318095
$ws.Range("T8").Value = 0.1890918412318094

# Row 9
$ws.Range("I9").Value = 0.8630927339690215
$ws.Range("J9").Value = 0.8630927339690215
$ws.Range("M9").Value = 25.27013633333333
$ws.Range("N9").Value = 75.81040899999999
$ws.Range("O9").Value = 0.2521761862513699
$ws.Range("P9").Value = 0.2521761862513699
$ws.Range("Q9").Value = 99.67897091512009
$ws.Range("R9").Value = 897.1107382360808
$ws.Range("S9").Value = 0.2176514340335761
$ws.Range("T9").Value = 0.2176514340335761

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.332886
$ws.Range("H10").Value = 0.998658
$ws.Range("I10").Value = 0.07283783531465635
$ws.Range("J10").Value = 0.07283783531465635
$ws.Range("M10").Value = 14.440165
$ws.Range("N10").Value = 43.320495
$ws.Range("O10").Value = 0.1441015470002482
$ws.Range("P10").Value = 0.1441015470002482
$ws.Range("Q10").Value = 4.80692876619
$ws.Range("R10").Value = 43.26235889571
$ws.Range("S10").Value = 0.01049604474899129
$ws.Range("T10").Value = 0.01049604474899129

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.332886
$ws.Range("H11").Value = 0.998658
$ws.Range("I11").Value = 0.07283783531465635
$ws.Range("J11").Value = 0.07283783531465635
$ws.Range("O11").Value = 0.3846359116098663
$ws.Range("P11").Value = 0.3846359116098662
$ws.Range("Q11").Value = 12.830656342808
$ws.Range("R11").Value = 115.475907085272
$ws.Range("S11").Value = 0.02801604718594215
$ws.Range("T11").Value = 0.02801604718594215

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.332886
$ws.Range("H12").Value = 0.998658
$ws.Range("I12").Value = 0.07283783531465635
$ws.Range("J12").Value = 0.07283783531465635
$ws.Range("M12").Value = 21.954262
$ws.Range("N12").Value = 65.862786
$ws.Range("O12").Value = 0.2190863551385157
$ws.Range("P12").Value = 0.2190863551385156
$ws.Range("Q12").Value = 7.308266460132001
$ws.Range("R12").Value = 65.774398141188
$ws.Range("S12").Value = 0.01595777585526752
$ws.Range("T12").Value = 0.01595777585526752

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.332886
$ws.Range("H13").Value = 0.998658
$ws.Range("I13").Value = 0.07283783531465635
$ws.Range("J13").Value = 0.07283783531465635
$ws.Range("M13").Value = 25.27013633333333
$ws.Range("N13").Value = 75.81040899999999
$ws.Range("O13").Value = 0.2521761862513699
$ws.Range("P13").Value = 0.2521761862513699
$ws.Range("Q13").Value = 8.412074603457999
$ws.Range("R13").Value = 75.708671431122
$ws.Range("S13").Value = 0.01836796752445539
$ws.Range("T13").Value = 0.01836796752445539
